# Update "想去人数" (column F) values on both the "展览" and "全部类型" sheets.
# Row 3 and row 13 are intentionally left unchanged.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 195
    4  = 263
    5  = 795
    6  = 250
    7  = 5981
    8  = 39
    9  = 68
    10 = 101
    11 = 51
    12 = 29
    14 = 181
    15 = 378
    16 = 29
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
